$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three data cells in row 2 (name, email, password) with the new
# values. The hyperlinks already attached to B2/C2 keep pointing at their
# original mailto: targets (khan39abcd@gmail.com / kishore@123) - only the
# displayed text changes.
$ws.Range("A2").Value = "Om ChaithanyaV"
$ws.Range("B2").Value = "omcv1989@gmail.com"
$ws.Range("C2").Value = "omcv902930"

# C2's hyperlink text no longer matches its target address, so Excel records
# an explicit display string ("kishore@123", the old password text) for it.
foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$C$2') {
        $h.TextToDisplay = "kishore@123"
    }
}

# Move the active selection to D2.
$ws.Range("D2").Select() | Out-Null
